$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1521.2222
$ws.Cells.Item(98, 9).Value = 1537.4231
$ws.Cells.Item(98, 10).Value = 1100
$ws.Cells.Item(98, 11).Value = 1537.4231
$ws.Cells.Item(98, 12).Value = 1100
$ws.Cells.Item(98, 13).Value = -39.42309999999998
$ws.Cells.Item(98, 14).Value = -4096

$ws.Cells.Item(122, 8).Value = 1521.2222
$ws.Cells.Item(122, 9).Value = 1537.4231
$ws.Cells.Item(122, 10).Value = 1100
$ws.Cells.Item(122, 11).Value = 4612.2693
$ws.Cells.Item(122, 12).Value = 3300
$ws.Cells.Item(122, 13).Value = -2162.2693
$ws.Cells.Item(122, 14).Value = -8200

$ws.Cells.Item(125, 8).Value = 2583.111
$ws.Cells.Item(125, 9).Value = 2500.3333
$ws.Cells.Item(125, 11).Value = 22502.9997
$ws.Cells.Item(125, 13).Value = -20042.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(30, 8).Value = 2690.125
$ws.Cells.Item(30, 10).Value = 2534.6667
$ws.Cells.Item(30, 12).Value = 2534.6667
$ws.Cells.Item(30, 14).Value = -2834.6667

$ws.Cells.Item(32, 8).Value = 125132
$ws.Cells.Item(32, 9).Value = 140148.38
$ws.Cells.Item(32, 10).Value = 5001
$ws.Cells.Item(32, 11).Value = 140148.38
$ws.Cells.Item(32, 12).Value = 5001
$ws.Cells.Item(32, 13).Value = -139861.38
$ws.Cells.Item(32, 14).Value = -5575

$ws.Cells.Item(61, 8).Value = 5454.9287
$ws.Cells.Item(61, 9).Value = 3685.5
$ws.Cells.Item(61, 10).Value = 9878.5
$ws.Cells.Item(61, 11).Value = 3685.5
$ws.Cells.Item(61, 12).Value = 9878.5
$ws.Cells.Item(61, 13).Value = -3473.5
$ws.Cells.Item(61, 14).Value = -10302.5

$ws.Cells.Item(74, 8).Value = 305809.94
$ws.Cells.Item(74, 9).Value = 667182.6
$ws.Cells.Item(74, 11).Value = 667182.6
$ws.Cells.Item(74, 13).Value = -666308.6

$ws.Cells.Item(77, 8).Value = 305809.94
$ws.Cells.Item(77, 9).Value = 667182.6
$ws.Cells.Item(77, 11).Value = 3335913
$ws.Cells.Item(77, 13).Value = -3331545

$ws.Cells.Item(110, 8).Value = 28790068
$ws.Cells.Item(110, 9).Value = 43183172
$ws.Cells.Item(110, 10).Value = 3859.9092
$ws.Cells.Item(110, 11).Value = 43183172
$ws.Cells.Item(110, 12).Value = 3859.9092
$ws.Cells.Item(110, 13).Value = -43181127
$ws.Cells.Item(110, 14).Value = -7949.9092

$ws.Cells.Item(126, 8).Value = 9999.5
$ws.Cells.Item(126, 9).Value = 9999.5
$ws.Cells.Item(126, 11).Value = 29998.5
$ws.Cells.Item(126, 13).Value = -27528.5

$ws.Cells.Item(132, 8).Value = 4510.05
$ws.Cells.Item(132, 9).Value = 3300.5334
$ws.Cells.Item(132, 11).Value = 9901.600199999999
$ws.Cells.Item(132, 13).Value = -7371.600199999999

$ws.Cells.Item(136, 8).Value = 5454.9287
$ws.Cells.Item(136, 9).Value = 3685.5
$ws.Cells.Item(136, 10).Value = 9878.5
$ws.Cells.Item(136, 11).Value = 11056.5
$ws.Cells.Item(136, 12).Value = 29635.5
$ws.Cells.Item(136, 13).Value = -8506.5
$ws.Cells.Item(136, 14).Value = -34735.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(102, 8).Value = 14837.818
$ws.Cells.Item(102, 9).Value = 4672.75
$ws.Cells.Item(102, 10).Value = 41944.668
$ws.Cells.Item(102, 11).Value = 4672.75
$ws.Cells.Item(102, 12).Value = 41944.668
$ws.Cells.Item(102, 13).Value = -1427.75
$ws.Cells.Item(102, 14).Value = -48434.668

$ws.Cells.Item(105, 8).Value = 76943640
$ws.Cells.Item(105, 9).Value = 76943640
$ws.Cells.Item(105, 11).Value = 76943640
$ws.Cells.Item(105, 13).Value = -76941893

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4757.5
$ws.Cells.Item(16, 9).Value = 4002.2
$ws.Cells.Item(16, 11).Value = 4002.2
$ws.Cells.Item(16, 13).Value = -3715.2

$ws.Cells.Item(99, 8).Value = 3270.75
$ws.Cells.Item(99, 9).Value = 3295.25
$ws.Cells.Item(99, 10).Value = 3246.25
$ws.Cells.Item(99, 11).Value = 3295.25
$ws.Cells.Item(99, 12).Value = 3246.25
$ws.Cells.Item(99, 13).Value = -1797.25
$ws.Cells.Item(99, 14).Value = -6242.25

$ws.Cells.Item(113, 8).Value = 4757.5
$ws.Cells.Item(113, 9).Value = 4002.2
$ws.Cells.Item(113, 11).Value = 4002.2
$ws.Cells.Item(113, 13).Value = -1832.2

$ws.Cells.Item(122, 8).Value = 92030.27
$ws.Cells.Item(122, 9).Value = 112259.336
$ws.Cells.Item(122, 11).Value = 336778.008
$ws.Cells.Item(122, 13).Value = -334328.008

$ws.Cells.Item(126, 8).Value = 3270.75
$ws.Cells.Item(126, 9).Value = 3295.25
$ws.Cells.Item(126, 10).Value = 3246.25
$ws.Cells.Item(126, 11).Value = 9885.75
$ws.Cells.Item(126, 12).Value = 9738.75
$ws.Cells.Item(126, 13).Value = -7415.75
$ws.Cells.Item(126, 14).Value = -14678.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 16393.21
$ws.Cells.Item(26, 9).Value = 132.48276
$ws.Cells.Item(26, 11).Value = 397.4482800000001
$ws.Cells.Item(26, 13).Value = -109.4482800000001

$ws.Cells.Item(61, 8).Value = 76.79412000000001
$ws.Cells.Item(61, 9).Value = 36.25
$ws.Cells.Item(61, 10).Value = 82.2
$ws.Cells.Item(61, 11).Value = 108.75
$ws.Cells.Item(61, 12).Value = 246.6
$ws.Cells.Item(61, 13).Value = 106.25
$ws.Cells.Item(61, 14).Value = -676.6

$ws.Cells.Item(95, 8).Value = 2000
$ws.Cells.Item(95, 9).Value = 2000
$ws.Cells.Item(95, 11).Value = 6000
$ws.Cells.Item(95, 13).Value = -3941

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2288.5454
$ws.Cells.Item(102, 9).Value = 2117.0625
$ws.Cells.Item(102, 11).Value = 2117.0625
$ws.Cells.Item(102, 13).Value = -495.0625

$ws.Cells.Item(122, 8).Value = 6606.788
$ws.Cells.Item(122, 9).Value = 8397.714
$ws.Cells.Item(122, 10).Value = 3472.6667
$ws.Cells.Item(122, 11).Value = 25193.142
$ws.Cells.Item(122, 12).Value = 10418.0001
$ws.Cells.Item(122, 13).Value = -22743.142
$ws.Cells.Item(122, 14).Value = -15318.0001

$ws.Cells.Item(132, 8).Value = 5062.3335
$ws.Cells.Item(132, 9).Value = 2873.3635
$ws.Cells.Item(132, 10).Value = 6567.25
$ws.Cells.Item(132, 11).Value = 8620.0905
$ws.Cells.Item(132, 12).Value = 19701.75
$ws.Cells.Item(132, 13).Value = -6090.0905
$ws.Cells.Item(132, 14).Value = -24761.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2782.8975
$ws.Cells.Item(22, 9).Value = 1889.9131
$ws.Cells.Item(22, 10).Value = 4066.5625
$ws.Cells.Item(22, 11).Value = 1889.9131
$ws.Cells.Item(22, 12).Value = 4066.5625
$ws.Cells.Item(22, 13).Value = -1594.9131
$ws.Cells.Item(22, 14).Value = -4656.5625

$ws.Cells.Item(27, 8).Value = 2782.8975
$ws.Cells.Item(27, 9).Value = 1889.9131
$ws.Cells.Item(27, 10).Value = 4066.5625
$ws.Cells.Item(27, 11).Value = 1889.9131
$ws.Cells.Item(27, 12).Value = 4066.5625
$ws.Cells.Item(27, 13).Value = -1782.9131
$ws.Cells.Item(27, 14).Value = -4280.5625

$ws.Cells.Item(40, 8).Value = 13696.6
$ws.Cells.Item(40, 9).Value = 13546.125
$ws.Cells.Item(40, 10).Value = 14298.5
$ws.Cells.Item(40, 11).Value = 13546.125
$ws.Cells.Item(40, 12).Value = 14298.5
$ws.Cells.Item(40, 13).Value = -13410.125
$ws.Cells.Item(40, 14).Value = -14570.5

$ws.Cells.Item(61, 8).Value = 2537.3333
$ws.Cells.Item(61, 9).Value = 1847.2858
$ws.Cells.Item(61, 11).Value = 1847.2858
$ws.Cells.Item(61, 13).Value = -1645.2858

$ws.Cells.Item(93, 8).Value = 6029.4614
$ws.Cells.Item(93, 9).Value = 7800
$ws.Cells.Item(93, 11).Value = 7800
$ws.Cells.Item(93, 13).Value = -6552

$ws.Cells.Item(113, 8).Value = 2537.3333
$ws.Cells.Item(113, 9).Value = 1847.2858
$ws.Cells.Item(113, 11).Value = 1847.2858
$ws.Cells.Item(113, 13).Value = 322.7141999999999

$ws.Cells.Item(127, 8).Value = 44357.5
$ws.Cells.Item(127, 10).Value = 44357.5
$ws.Cells.Item(127, 12).Value = 44357.5
$ws.Cells.Item(127, 14).Value = -54277.5

$ws.Cells.Item(132, 8).Value = 5995.2
$ws.Cells.Item(132, 9).Value = 4815.933
$ws.Cells.Item(132, 10).Value = 9533
$ws.Cells.Item(132, 11).Value = 14447.799
$ws.Cells.Item(132, 12).Value = 28599
$ws.Cells.Item(132, 13).Value = -11917.799
$ws.Cells.Item(132, 14).Value = -33659

$ws.Cells.Item(136, 8).Value = 4321.65
$ws.Cells.Item(136, 9).Value = 3656.2
$ws.Cells.Item(136, 10).Value = 4720.92
$ws.Cells.Item(136, 11).Value = 10968.6
$ws.Cells.Item(136, 12).Value = 14162.76
$ws.Cells.Item(136, 13).Value = -8418.599999999999
$ws.Cells.Item(136, 14).Value = -19262.76

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value = 11302.333
$ws.Cells.Item(52, 9).Value = 3653.6667
$ws.Cells.Item(52, 10).Value = 26599.666
$ws.Cells.Item(52, 11).Value = 3653.6667
$ws.Cells.Item(52, 12).Value = 26599.666
$ws.Cells.Item(52, 13).Value = -3427.6667
$ws.Cells.Item(52, 14).Value = -27051.666

$ws.Cells.Item(61, 8).Value = 14755.728
$ws.Cells.Item(61, 9).Value = 11851.143
$ws.Cells.Item(61, 11).Value = 11851.143
$ws.Cells.Item(61, 13).Value = -11559.143

$ws.Cells.Item(96, 8).Value = 67668.336
$ws.Cells.Item(96, 9).Value = 3000
$ws.Cells.Item(96, 10).Value = 100002.5
$ws.Cells.Item(96, 11).Value = 3000
$ws.Cells.Item(96, 12).Value = 100002.5
$ws.Cells.Item(96, 13).Value = -1627
$ws.Cells.Item(96, 14).Value = -102748.5

$ws.Cells.Item(107, 8).Value = 2053.4546
$ws.Cells.Item(107, 9).Value = 2101.9375
$ws.Cells.Item(107, 11).Value = 6305.8125
$ws.Cells.Item(107, 13).Value = -4385.8125

$ws.Cells.Item(136, 8).Value = 5828.0713
$ws.Cells.Item(136, 9).Value = 2182.3333
$ws.Cells.Item(136, 11).Value = 6546.999899999999
$ws.Cells.Item(136, 13).Value = -3996.999899999999
